# Auto-generated script to apply numeric value updates to Behemoth_Profits workbook
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3475.3901
$ws.Range("I15").Value = 3475.3901
$ws.Range("K15").Value = 10426.1703
$ws.Range("M15").Value = -10257.1703
$ws.Range("H28").Value = 499.33334
$ws.Range("I28").Value = 399.2
$ws.Range("K28").Value = 399.2
$ws.Range("M28").Value = 85.80000000000001
$ws.Range("H31").Value = 167566.5
$ws.Range("I31").Value = 299.75
$ws.Range("J31").Value = 502100
$ws.Range("K31").Value = 899.25
$ws.Range("L31").Value = 1506300
$ws.Range("M31").Value = -669.25
$ws.Range("N31").Value = -1506760
$ws.Range("H55").Value = 629.8
$ws.Range("J55").Value = 866.3333
$ws.Range("L55").Value = 866.3333
$ws.Range("N55").Value = -1294.3333
$ws.Range("H58").Value = 15773.818
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 15773.818
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 47321.454
$ws.Range("N58").Value = -47621.454
$ws.Range("H80").Value = 770.3125
$ws.Range("I80").Value = 543.1429000000001
$ws.Range("J80").Value = 947
$ws.Range("K80").Value = 1629.4287
$ws.Range("L80").Value = 2841
$ws.Range("M80").Value = -631.4287000000002
$ws.Range("N80").Value = -4837
$ws.Range("H83").Value = 770.3125
$ws.Range("I83").Value = 543.1429000000001
$ws.Range("J83").Value = 947
$ws.Range("K83").Value = 4888.2861
$ws.Range("L83").Value = 8523
$ws.Range("M83").Value = 103.7138999999997
$ws.Range("N83").Value = -18507
$ws.Range("H111").Value = 11516
$ws.Range("I111").Value = 3000
$ws.Range("K111").Value = 9000
$ws.Range("M111").Value = -5933
$ws.Range("H113").Value = 41677584
$ws.Range("I113").Value = 11122332
$ws.Range("K113").Value = 11122332
$ws.Range("M113").Value = -11119078
$ws.Range("H138").Value = 3696.3972
$ws.Range("J138").Value = 4002.6897
$ws.Range("L138").Value = 12008.0691
$ws.Range("N138").Value = -22288.0691

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8785486
$ws.Range("I32").Value = 10212449
$ws.Range("K32").Value = 10212449
$ws.Range("M32").Value = -10212162
$ws.Range("H122").Value = 2774.4443
$ws.Range("I122").Value = 1331.3158
$ws.Range("K122").Value = 3993.9474
$ws.Range("M122").Value = -1543.9474
$ws.Range("H129").Value = 78332
$ws.Range("I129").Value = 59998
$ws.Range("J129").Value = 115000
$ws.Range("K129").Value = 59998
$ws.Range("L129").Value = 115000
$ws.Range("M129").Value = -54998
$ws.Range("N129").Value = -125000
$ws.Range("H132").Value = 5863.1816
$ws.Range("I132").Value = 2743.476
$ws.Range("J132").Value = 15942.23
$ws.Range("K132").Value = 8230.428
$ws.Range("L132").Value = 47826.69
$ws.Range("M132").Value = -5700.428
$ws.Range("N132").Value = -52886.69

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2037.7
$ws.Range("I107").Value = 1668.8572
$ws.Range("J107").Value = 2898.3333
$ws.Range("K107").Value = 1668.8572
$ws.Range("L107").Value = 2898.3333
$ws.Range("M107").Value = 251.1428000000001
$ws.Range("N107").Value = -6738.3333
$ws.Range("H134").Value = 27547.365
$ws.Range("I134").Value = 3277.8718
$ws.Range("K134").Value = 9833.615399999999
$ws.Range("M134").Value = -7298.615399999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2517.4412
$ws.Range("I132").Value = 2362.2812
$ws.Range("K132").Value = 7086.8436
$ws.Range("M132").Value = -4556.8436
$ws.Range("H134").Value = 478715.94
$ws.Range("I134").Value = 668689.2
$ws.Range("K134").Value = 2006067.6
$ws.Range("M134").Value = -2003532.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 715684
$ws.Range("J92").Value = 798
$ws.Range("L92").Value = 2394
$ws.Range("N92").Value = -4890
$ws.Range("H122").Value = 595.53845
$ws.Range("I122").Value = 536.8
$ws.Range("J122").Value = 632.25
$ws.Range("K122").Value = 4831.2
$ws.Range("L122").Value = 5690.25
$ws.Range("M122").Value = -2381.2
$ws.Range("N122").Value = -10590.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 24999.666
$ws.Range("J40").Value = 24999.666
$ws.Range("L40").Value = 24999.666
$ws.Range("N40").Value = -25301.666
$ws.Range("H80").Value = 3996.5
$ws.Range("J80").Value = 3996.5
$ws.Range("L80").Value = 3996.5
$ws.Range("N80").Value = -5992.5
$ws.Range("H83").Value = 3996.5
$ws.Range("J83").Value = 3996.5
$ws.Range("L83").Value = 19982.5
$ws.Range("N83").Value = -29966.5
$ws.Range("H102").Value = 2869.9429
$ws.Range("I102").Value = 2400.68
$ws.Range("K102").Value = 2400.68
$ws.Range("M102").Value = -778.6799999999998
$ws.Range("H122").Value = 3715.3333
$ws.Range("I122").Value = 4016.6365
$ws.Range("K122").Value = 12049.9095
$ws.Range("M122").Value = -9599.9095
$ws.Range("H126").Value = 4049.4285
$ws.Range("I126").Value = 3200.4
$ws.Range("K126").Value = 9601.200000000001
$ws.Range("M126").Value = -7131.200000000001
$ws.Range("H132").Value = 22224906
$ws.Range("I132").Value = 25643708
$ws.Range("J132").Value = 2699
$ws.Range("K132").Value = 76931124
$ws.Range("L132").Value = 8097
$ws.Range("M132").Value = -76928594
$ws.Range("N132").Value = -13157

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 148217.28
$ws.Range("I7").Value = 3750.5
$ws.Range("K7").Value = 3750.5
$ws.Range("M7").Value = -3638.5
$ws.Range("H46").Value = 2448.125
$ws.Range("I46").Value = 2296.25
$ws.Range("J46").Value = 2600
$ws.Range("K46").Value = 2296.25
$ws.Range("L46").Value = 2600
$ws.Range("M46").Value = -2108.25
$ws.Range("N46").Value = -2976
$ws.Range("H126").Value = 148217.28
$ws.Range("I126").Value = 3750.5
$ws.Range("K126").Value = 11251.5
$ws.Range("M126").Value = -8781.5
$ws.Range("H132").Value = 837904.0600000001
$ws.Range("I132").Value = 1254656.2
$ws.Range("J132").Value = 4399.75
$ws.Range("K132").Value = 3763968.6
$ws.Range("L132").Value = 13199.25
$ws.Range("M132").Value = -3761438.6
$ws.Range("N132").Value = -18259.25
$ws.Range("H136").Value = 354917.34
$ws.Range("I136").Value = 337499.66
$ws.Range("K136").Value = 1012498.98
$ws.Range("M136").Value = -1009948.98

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 11324
$ws.Range("J74").Value = 11324
$ws.Range("L74").Value = 11324
$ws.Range("N74").Value = -13196
$ws.Range("H77").Value = 11324
$ws.Range("J77").Value = 11324
$ws.Range("L77").Value = 33972
$ws.Range("N77").Value = -43332
$ws.Range("H113").Value = 843.7742
$ws.Range("I113").Value = 771.95
$ws.Range("J113").Value = 974.36365
$ws.Range("K113").Value = 2315.85
$ws.Range("L113").Value = 2923.09095
$ws.Range("M113").Value = -145.8500000000004
$ws.Range("N113").Value = -7263.09095
$ws.Range("H122").Value = 5741.0625
$ws.Range("I122").Value = 3912.5386
$ws.Range("J122").Value = 13664.667
$ws.Range("K122").Value = 11737.6158
$ws.Range("L122").Value = 40994.001
$ws.Range("M122").Value = -9287.6158
$ws.Range("N122").Value = -45894.001
$ws.Range("H126").Value = 3713.4707
$ws.Range("I126").Value = 2785.5588
$ws.Range("K126").Value = 8356.6764
$ws.Range("M126").Value = -5886.6764

# --- Special case: ALC!M58 is fully cleared (cell removed) per diff ---
$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("M58").ClearContents()
